$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.08707909936379449
$ws.Range("C2").Value = 0.1849415347999282
$ws.Range("D2").Value = 0.6087619768527019
$ws.Range("B3").Value = 0.2343920533817323
$ws.Range("C3").Value = 0.5621024059030257
$ws.Range("B4").Value = 0.7136637901383804
$ws.Range("B5").Value = 0.6611739089257214
$ws.Range("C5").Value = 0.2424661772861064
$ws.Range("D5").Value = 0.343437461374561
$ws.Range("E5").Value = 0.2316648163334098
$ws.Range("B6").Value = 0.52173245889808
$ws.Range("C6").Value = 0.3916150992080441
$ws.Range("D6").Value = 0.0731403976960974
$ws.Range("B7").Value = 0.6205335770846683
$ws.Range("C7").Value = 0.1160860902143776
$ws.Range("B8").Value = 0.3924296930042493
$ws.Range("B9").Value = 0.4747974866332783
$ws.Range("C9").Value = -0.06986322861814367
$ws.Range("D9").Value = 0.2092250849631835
$ws.Range("E9").Value = 0.02501891626858249
$ws.Range("B10").Value = 0.2191407266580338
$ws.Range("C10").Value = 0.2222972325179511
$ws.Range("D10").Value = -0.09185372756793091
$ws.Range("B11").Value = 0.224871884804233
$ws.Range("C11").Value = -0.0167560416292489
$ws.Range("B12").Value = 0.1712949347374914
$ws.Range("B13").Value = -0.01257980384821508
$ws.Range("C13").Value = 0.1395187281562265
$ws.Range("D13").Value = 0.1036987463175434
$ws.Range("B14").Value = 0.2732997081603292
$ws.Range("C14").Value = 0.04948772950879379
$ws.Range("B15").Value = 0.1075450023515467
